# Trade #20 closed at 2026-02-17 15:18:43 - unknown UNKNOWN +0.000%
# Applies the new trade row to "All Trades" and "MarketMaking" sheets,
# and updates the rolled-up summary figures on "Summary" and
# "Strategy Status" accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.78
$summary.Range("B4").Value = -0.23
$summary.Range("B5").Value = -0.23
$summary.Range("B6").Value = 20
$summary.Range("B8").Value = 10
$summary.Range("B9").Value = 25

# ---------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.78
$status.Range("D4").Value = 20
$status.Range("E4").Value = -0.23
$status.Range("F4").Value = -0.22
$status.Range("G4").Value = 25

# ---------------------------------------------------------------
# 3) Append new trade row (#20) to both "All Trades" and
#    "MarketMaking" sheets at row 21.
# ---------------------------------------------------------------
$newRow = @{
    A = 20
    B = "'2026-02-17"
    C = "15:18:36"
    D = "MarketMaking"
    E = "UP"
    F = 0.92
    G = 0.9
    H = "CLOSED"
    I = -2.1739
    J = -0.02
    K = 99.78
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.15
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A21").Value = $newRow.A
    $ws.Range("B21").Value = $newRow.B
    $ws.Range("C21").Value = $newRow.C
    $ws.Range("D21").Value = $newRow.D
    $ws.Range("E21").Value = $newRow.E
    $ws.Range("F21").Value = $newRow.F
    $ws.Range("G21").Value = $newRow.G
    $ws.Range("H21").Value = $newRow.H
    $ws.Range("I21").Value = $newRow.I
    $ws.Range("J21").Value = $newRow.J
    $ws.Range("K21").Value = $newRow.K
    $ws.Range("L21").Value = $newRow.L
    $ws.Range("M21").Value = $newRow.M
    $ws.Range("N21").Value = $newRow.N
    $ws.Range("O21").Value = $newRow.O
    $ws.Range("P21").Value = $newRow.P
    $ws.Range("Q21").Value = $newRow.Q
}
